$d = $word.ActiveDocument

# 1. "lead" -> "led" in "Our research lead us to select PYQT..."
$d.Content.Find.Execute("research lead us", $true, $false, $false, $false, $false, $true, 1, $false, "research led us", 2)

# 2. "or" -> "and" in "Software can be written using both a Desktop or a Web Application approach."
$d.Content.Find.Execute("Desktop or a Web", $true, $false, $false, $false, $false, $true, 1, $false, "Desktop and a Web", 2)

# 3. "final outcome" -> "outcome" in "...have flow until the final outcome was reached."
$d.Content.Find.Execute("the final outcome was reached", $true, $false, $false, $false, $false, $true, 1, $false, "the outcome was reached", 2)

# 4. "step by step" -> "step-by-step" in "...that gives a step by step guide."
$d.Content.Find.Execute("gives a step by step guide", $true, $false, $false, $false, $false, $true, 1, $false, "gives a step-by-step guide", 2)
